# Blackpool To-Do List update:
#  - Add a new (hidden, filtered-out) task row at row 9:
#       Task:   Add "Admin" Menu to Blackpool with BASECAP options, AIRBOSS start/stop, more?
#       Priority: 1-High
#       Status: Complete
#       % Complete: 100%
#  - Move the active selection from C15 to B10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 'Add "Admin" Menu to Blackpool with BASECAP options, AIRBOSS start/stop, more?'
$ws.Range("C9").Value = "1-High"
$ws.Range("D9").Value = "Complete"
$ws.Range("E9").Value = 1

# Row is filtered out by the table's autofilter (Status <> "Not Started"/blank),
# so it ends up hidden just like the other non-"Not Started" rows (4, 5, 8).
$ws.Rows.Item(9).Hidden = $true

# Update the active selection/cell
$ws.Range("B10").Select()
